$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.519.38"
$ws.Range("E2").Value = "  +5.67%  "
$ws.Range("D3").Value = "1.708.53"
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("E4").Value = "  -0.06%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "222.38"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.29%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("E9").Value = "  +3.36%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0648"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +6.52%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0911"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Value = "1.955.09"
$ws.Range("E12").Value = "  +4.30%  "
$ws.Range("D13").Value = "1.702.83"
$ws.Range("E13").Value = "  +3.96%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.613"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +4.06%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "10.16"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +8.18%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "4.20"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +8.44%  "
$ws.Range("D17").Value = "31.505.54"
$ws.Range("E17").Value = "  +5.58%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "67.30"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +4.72%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "251.02"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +5.09%  "
$ws.Range("D20").Value = "0.0₃0724"
$ws.Range("E20").Value = "  +3.17%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("E23").Value = "  +3.00%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.17"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.66%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "159.62"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "16.05"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +3.46%  "
$ws.Range("E27").Value = "  +3.32%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "6.79"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("E29").Value = "  -0.12%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.90"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +15.46%  "
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +4.10%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.40"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +6.53%  "
$ws.Range("D34").Value = "1.526.93"
$ws.Range("E34").Value = "  +7.50%  "
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.613"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +9.12%  "
$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "82.53"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +8.27%  "
$ws.Range("E39").Value = "  +4.63%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.71"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.04"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +5.11%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.853"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.66%  "
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("E46").Value = "  -0.05%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "52.19"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +5.91%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "5.60"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +4.62%  "
$ws.Range("E49").Value = "  +3.72%  "
$ws.Range("E50").Value = "  +9.84%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "93.67"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.52%  "
